# Updates the cryptocurrency price/volume table on Sheet1 (rows 2-51) to
# reflect the latest scrape, matching the "Updated cryptos list ... with
# GitHub Actions" commit. Most rows keep the same coin but refresh the
# Price (D) and Volume(1h) (E) columns; a few rows also swap which coin
# occupies that rank (columns B/C change together with D/E).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price values in column D are text (e.g. "45.930.12", "1.00", "2.80") and
# must NOT be auto-converted to numbers by Excel (which would both change
# the cell type and silently drop meaningful trailing zeros / thousands
# grouping). Temporarily force a Text number format while assigning the
# value, then restore the default "Normal" style so no residual
# formatting is left behind on the cell.
function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "45.930.12"
$ws.Range("E2").Value = "  -1.24%  "

Set-TextValue $ws.Range("D3") "2.618.25"
$ws.Range("E3").Value = "  -0.10%  "

Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  -0.04%  "

Set-TextValue $ws.Range("D5") "309.57"
$ws.Range("E5").Value = "  -1.45%  "

Set-TextValue $ws.Range("D6") "98.61"
$ws.Range("E6").Value = "  -3.49%  "

Set-TextValue $ws.Range("D7") "0.596"
$ws.Range("E7").Value = "  -1.15%  "

$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("E9").Value = "  -2.16%  "

Set-TextValue $ws.Range("D10") "38.78"
$ws.Range("E10").Value = "  -0.62%  "

Set-TextValue $ws.Range("D11") "54.11"
$ws.Range("E11").Value = "  -0.93%  "

$ws.Range("E12").Value = "  -0.39%  "

Set-TextValue $ws.Range("D13") "8.05"
$ws.Range("E13").Value = "  -3.71%  "

Set-TextValue $ws.Range("D14") "3.017.85"
$ws.Range("E14").Value = "  -0.10%  "

$ws.Range("E15").Value = "  +0.59%  "

Set-TextValue $ws.Range("D16") "2.624.93"
$ws.Range("E16").Value = "  +0.00%  "

Set-TextValue $ws.Range("D17") "0.916"
$ws.Range("E17").Value = "  -0.17%  "

Set-TextValue $ws.Range("D18") "14.83"
$ws.Range("E18").Value = "  -2.16%  "

Set-TextValue $ws.Range("D19") "45.938.35"
$ws.Range("E19").Value = "  -1.78%  "

Set-TextValue $ws.Range("D20") "0.0000102"
$ws.Range("E20").Value = "  -1.29%  "

Set-TextValue $ws.Range("D21") "6.75"
$ws.Range("E21").Value = "  -0.54%  "

Set-TextValue $ws.Range("D22") "12.71"
$ws.Range("E22").Value = "  -4.91%  "

Set-TextValue $ws.Range("D23") "74.74"
$ws.Range("E23").Value = "  +5.19%  "

Set-TextValue $ws.Range("D24") "282.12"
$ws.Range("E24").Value = "  +10.36%  "

Set-TextValue $ws.Range("D25") "3.03"
$ws.Range("E25").Value = "  -2.67%  "

Set-TextValue $ws.Range("D26") "2.26"
$ws.Range("E26").Value = "  +1.90%  "

Set-TextValue $ws.Range("D27") "29.63"
$ws.Range("E27").Value = "  +4.58%  "

$ws.Range("E28").Value = "  -0.18%  "

$ws.Range("E29").Value = "  +0.13%  "

Set-TextValue $ws.Range("D30") "10.45"
$ws.Range("E30").Value = "  -2.42%  "

Set-TextValue $ws.Range("D31") "38.71"
$ws.Range("E31").Value = "  -7.16%  "

$ws.Range("E32").Value = "  -3.43%  "

Set-TextValue $ws.Range("D33") "6.23"
$ws.Range("E33").Value = "  -0.02%  "

$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D34") "2.29"
$ws.Range("E34").Value = "  +0.05%  "

$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws.Range("D35") "3.60"
$ws.Range("E35").Value = "  -4.34%  "

Set-TextValue $ws.Range("D36") "156.88"
$ws.Range("E36").Value = "  +1.87%  "

Set-TextValue $ws.Range("D37") "0.0836"
$ws.Range("E37").Value = "  -0.47%  "

Set-TextValue $ws.Range("D38") "2.80"
$ws.Range("E38").Value = "  -2.35%  "

Set-TextValue $ws.Range("D39") "0.123"
$ws.Range("E39").Value = "  +3.10%  "

Set-TextValue $ws.Range("D40") "0.124"
$ws.Range("E40").Value = "  +0.50%  "

Set-TextValue $ws.Range("D41") "15.82"
$ws.Range("E41").Value = "  -7.61%  "

Set-TextValue $ws.Range("D42") "22.18"
$ws.Range("E42").Value = "  +4.80%  "

Set-TextValue $ws.Range("D43") "0.0327"
$ws.Range("E43").Value = "  -1.01%  "

$ws.Range("B44").Value = "NEARProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D44") "3.53"
$ws.Range("E44").Value = "  -3.31%  "

$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D45") "4.02"
$ws.Range("E45").Value = "  -6.90%  "

Set-TextValue $ws.Range("D46") "2.109.08"
$ws.Range("E46").Value = "  +3.42%  "

Set-TextValue $ws.Range("D47") "0.999"
$ws.Range("E47").Value = "  +0.08%  "

Set-TextValue $ws.Range("D48") "93.89"
$ws.Range("E48").Value = "  +2.55%  "

$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D49") "109.89"
$ws.Range("E49").Value = "  -3.06%  "

$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D50") "9.15"
$ws.Range("E50").Value = "  -1.18%  "

Set-TextValue $ws.Range("D51") "2.870.58"
$ws.Range("E51").Value = "  -0.37%  "
